# The edit shuffles the observation records currently sitting in rows
# 13-21 of the active sheet: each row's full set of field values moves to a
# different row among that same block (a permutation), with no rows added
# or removed. Concretely (new row <- source row):
#   13<-15  14<-17  15<-13  16<-20  17<-16  18<-14  19<-21  20<-18  21<-19
#
# Because several rows both give AND receive data (it's a set of cycles,
# not a simple swap), we must snapshot every source cell's value BEFORE
# writing anything, then write the new rows from that snapshot. Otherwise
# an early write could clobber a value that's still needed as a source for
# a later write.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that carry data in this block of rows.
$cols = @("A","B","C","D","E","F","G","H","I","K","L","M","N","P","Q","R","S","T","U","V","W","Y","AA","AD","AE","AG","AT","AW","AX","AY")

# Columns whose text values look like numbers/dates ("1", "2023-09-14", ...)
# in the source data, even though the underlying field is really text. Left
# to its normal smart-entry behaviour, Excel would silently reinterpret
# those as a real number / date serial when assigned through .Value, so for
# these columns we force text by leading the string with an apostrophe
# (standard Excel "treat as text" quoting) instead of auto-detecting type.
$forceTextCols = @("I","Y","AA")

# new row number -> row number currently holding the data that should end
# up there.
$rowMap = @{
    13 = 15
    14 = 17
    15 = 13
    16 = 20
    17 = 16
    18 = 14
    19 = 21
    20 = 18
    21 = 19
}

# 1) Snapshot current contents (Value2 preserves number/text/bool typing)
# of every row involved, for every used column.
$snapshot = @{}
foreach ($srcRow in 13..21) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$srcRow").Value2
    }
    $snapshot[$srcRow] = $rowVals
}

# 2) Write the snapshot back out in the new row order.
foreach ($newRow in 13..21) {
    $srcRow = $rowMap[$newRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $val = $rowVals[$c]
        if ($forceTextCols -contains $c -and $val -ne $null -and $val -ne "") {
            $ws.Range("$c$newRow").Value = "'" + $val
        } else {
            $ws.Range("$c$newRow").Value = $val
        }
    }
}
